# envios y arreglo del serviciosInternos movile y checkbox
#
# Appends two new survey-response rows (3 and 4) to the "Respuestas" sheet,
# mirroring the existing row 2 for the repeated contact (Luciano Albani,
# Mar del Plata) but with their own phone/budget-range/start-time/timestamp
# answers, and makes sure the sheet view is explicitly left-to-right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet should read left-to-right (matches the source survey export).
$ws.DisplayRightToLeft = $false

# --- Row 3 -----------------------------------------------------------
$ws.Range("A3").Value = "Luciano"
$ws.Range("B3").Value = "Albani"
$ws.Range("C3").Value = "Mar del Plata"
$ws.Range("D3").Value = "luchoalbanix1@gmail.com"
$ws.Range("E3").Value = "+54 2234480301"
$ws.Range("F3").Value = "rango4"
$ws.Range("G3").Value = "En un mes"
$ws.Range("H3").Value = "31/5/2025, 1:49:42 a.m."

# --- Row 4 -----------------------------------------------------------
$ws.Range("A4").Value = "Luciano"
$ws.Range("B4").Value = "Albani"
$ws.Range("C4").Value = "Mar del Plata"
$ws.Range("D4").Value = "luchoalbanix1@gmail.com"
$ws.Range("E4").Value = "+1242 2234480301"
$ws.Range("F4").Value = "rango5"
$ws.Range("G4").Value = "Inmediatamente"
$ws.Range("H4").Value = "31/5/2025, 3:15:56 a.m."

# The "number stored as text" ignored-error banner should keep covering
# the whole used range (header + all response rows) now that it grew.
$rng = $ws.Range("A1:H4")
$rng.Errors.Item(1).Ignore = $true
